$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update cell values (rows 2-14) to their new content.
# ---------------------------------------------------------------------------

# Row 2 - Associate Professor / 2019 / Universidad El Bosque (hyperlink) / Bogota
$ws.Range("A2").Value = "Associate Professor"
$ws.Range("C2").Value = "\href{https://www.unbosque.edu.co/}{Universidad El Bosque}"
$ws.Range("D2").Value = "Bogota, Colombia"
$ws.Range("E2").Value = "Quantitative Methods II (Psychology MSc) (2019)"

# Row 3 - Assistant Professor group (first row keeps who/where, rest is cleared)
$ws.Range("A3").Value = "Assistant Professor"
$ws.Range("B3").Value = "2015 - 2018"
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = $null
$ws.Range("E3").Value = "Quantitative Methods II (Psychology MSc) (2017-2018)"

# Row 4 - Assistant Professor / Research Degree Project (2018)
$ws.Range("A4").Value = "Assistant Professor"
$ws.Range("B4").Value = $null
$ws.Range("E4").Value = "Research Degree Project (2018)"

# Row 5 - Assistant Professor / Quantitative Methods I (2017)
$ws.Range("A5").Value = "Assistant Professor"
$ws.Range("B5").Value = $null
$ws.Range("E5").Value = "Quantitative Methods I (Psychology MSc) (2017)"

# Row 6 - Assistant Professor / Sources and Documentation Styles (2015)
$ws.Range("A6").Value = "Assistant Professor"
$ws.Range("B6").Value = $null
$ws.Range("E6").Value = "Sources and Documentation Styles in Psychology (2015)"

# Row 7 - Cathedratic Professor group
$ws.Range("A7").Value = "Cathedratic Professor"
$ws.Range("B7").Value = "2015 - 2016"
$ws.Range("C7").Value = "\href{https://www.unisabana.edu.co/}{Universidad de La Sabana}"
$ws.Range("D7").Value = "Chia, Colombia"
$ws.Range("E7").Value = "Evolution and Development of Vocal Communication: Songs, Fashion, and Language (2016)"

# Row 8 - Cathedratic Professor / Inferential Statistics
$ws.Range("A8").Value = "Cathedratic Professor"
$ws.Range("B8").Value = $null
$ws.Range("E8").Value = "Inferential Statistics (2015 - 2016)"

# Row 9 - Cathedratic Professor / Descriptive Statistics
$ws.Range("A9").Value = "Cathedratic Professor"
$ws.Range("B9").Value = $null
$ws.Range("E9").Value = "Descriptive Statistics (2015 - 2016)"

# Row 10 - Teaching Assistant group
$ws.Range("A10").Value = "Teaching Assistant"
$ws.Range("B10").Value = "2012 - 2014"
$ws.Range("C10").Value = "\href{https://www.stir.ac.uk/}{University of Stirling}"
$ws.Range("D10").Value = "Stirling, UK"
$ws.Range("E10").Value = "Animal Behaviour (lecture on vocal communication) (2012)"

# Row 11 - Teaching Assistant / Quantitative Methods (several lectures ...)
$ws.Range("A11").Value = "Teaching Assistant"
$ws.Range("B11").Value = $null
$ws.Range("E11").Value = "Quantitative Methods (Psychology MSc – several lectures, practical supervision, one-on-one teaching) (2012-2014)"

# Row 12 - Teaching Assistant / Cognition Module
$ws.Range("A12").Value = "Teaching Assistant"
$ws.Range("B12").Value = $null
$ws.Range("E12").Value = "Cognition Module (leading research projects in psychoacoustics) (2012-2014)"

# Row 13 - Auxiliar Professor group
$ws.Range("A13").Value = "Auxiliar Professor"
$ws.Range("B13").Value = 2010
$ws.Range("C13").Value = "\href{https://www.upn.edu.co/}{Universidad Pedagógica Nacional}"
$ws.Range("D13").Value = "Bogota, Colombia"
$ws.Range("E13").Value = "Research Project I (2010)"

# Row 14 - Auxiliar Professor / Research Lab II
$ws.Range("A14").Value = "Auxiliar Professor"
$ws.Range("E14").Value = "Research Lab II (2010)"

# ---------------------------------------------------------------------------
# 2. Add 7 blank (but formatted) rows below the data, rows 15-21.
# ---------------------------------------------------------------------------
$blank = $ws.Range("A15:E21")
$blank.Value = ""

# ---------------------------------------------------------------------------
# 3. Re-format the whole table: left/top aligned, wrapped text.
# ---------------------------------------------------------------------------
$full = $ws.Range("A1:E21")
$full.HorizontalAlignment = -4131
$full.VerticalAlignment = -4160
$full.WrapText = $true

# ---------------------------------------------------------------------------
# 4. Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.28515625
$ws.Columns.Item(2).ColumnWidth = 10.5703125
$ws.Columns.Item(3).ColumnWidth = 62.28515625
$ws.Columns.Item(4).ColumnWidth = 16.7109375
$ws.Columns.Item(5).ColumnWidth = 84.5703125

$ws.StandardWidth = 85.5703125

# ---------------------------------------------------------------------------
# 5. Row height: row 11 wraps to two lines.
# ---------------------------------------------------------------------------
$ws.Rows.Item(11).RowHeight = 30

# ---------------------------------------------------------------------------
# 6. Selection.
# ---------------------------------------------------------------------------
$ws.Range("C27").Select()

Write-Host "edit applied"
